# Diario.docx – add diary entries 05/06/2018 .. 08/06/2018
#
# The existing last real entry ("Añadir grupos invitados al crear
# usuario") carries a hidden _GoBack bookmark between the "invitados "
# and "al crear usuario" runs (left behind by Word at the point the
# author last typed). The new commit removes that bookmark (since more
# text follows now), appends several new Fechas/Entrada paragraph
# pairs, and leaves a fresh _GoBack bookmark collapsed at the very end
# of the final new paragraph. The old trailing empty paragraph (an
# artifact paragraph with numPr ilvl=0/numId=0) is dropped.

$d = $word.ActiveDocument

# --- 1. Drop the stray _GoBack bookmark -----------------------------
# Removing it (instead of rewriting the paragraph's text) keeps the
# "Añadir grupos " / "invitados " / "al crear usuario" runs intact and
# separate, exactly as they were except for the bookmark.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2. Build the OOXML for the new paragraphs -----------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml = @'
<w:p><w:pPr><w:pStyle w:val="Fechas"/></w:pPr><w:r><w:t>05/06/2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t>Lista de usuarios del grupo</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Fechas"/></w:pPr><w:r><w:t>06/06/2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t xml:space="preserve">Separar la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>activity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con un fragmento y botón de tipos  de turnos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Fechas"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>07/06/2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t>Borrar usuario marcándolo como no activo</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t xml:space="preserve">Inicializar activo en las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>functions</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Fechas"/></w:pPr><w:r><w:t>08/06/2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t xml:space="preserve">Añadido </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShiftTypeFragment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>recycler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>view</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de tipos de turno</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$newParagraphsXml = $newParagraphsXml.Trim() -replace '<w:p>', ('<w:p ' + $wNs + '>')

# --- 3. Insert the new paragraphs in place of the trailing empty one -
# The document's final paragraph is an empty artifact paragraph
# (style Entrada, numPr ilvl=0/numId=0). Collapsing a range to its
# start and inserting the OOXML fragment there replaces that throwaway
# paragraph with the real new content.
$lastParagraph = $d.Paragraphs.Last
$insertionPoint = $lastParagraph.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertXML($newParagraphsXml)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
